$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = " * On September 4, 2007, a murder took place at the residence of Haji Muhammad, father of the complainant, Muhammad Arshad. * The complainant alleged that Muhammad Makki, Muhammad Mazhar, Muhammad Akhtar, Ghulam Mustafa, Muhammad Aslam, and Ghulam Sarwar, along with Abdul Ghafoor, committed the murder. * The accused allegedly used firearms, causing injuries to Muhammad Siddique and Muhammad Umair, and also caused the death of Haji Muhammad, Muhammad Umar, and Sajida Bibi. * The accused were charged with murder, rioting, and terrorism, and were convicted and sentenced by the Anti-Terrorism Court-I, Multan.  "

$ws.Range("B2").Value = " * Muhammad Arshad, the complainant, testified that he and his brothers were sleeping on the roof of their residential room when they heard the sound of motorcycles and saw the accused making firing upon the deceased and injured. * Muhammad Arshad specifically nominated each accused with specific firearms and attributed specific firearm injuries to each accused. * Muhammad Akhtar, one of the accused, testified that he was asleep at the time of the incident and did not participate in the crime. * Muhammad Umair, one of the injured witnesses, was not produced as a prosecution witness, and it was presumed that he would not have supported the prosecution version if he had been produced.  "

$ws.Range("C2").Value = " * The Supreme Court of Pakistan acquitted Muhammad Akhtar, one of the accused, and dismissed the appeals filed by the complainant against the acquittal of the co-accused and for enhancement of the sentence of the convicts. * The court held that the prosecution failed to prove the charge against Muhammad Akhtar beyond reasonable doubt due to the lack of independent corroboration of the recovery of a 44 bore rifle from him, and the contradictions and dishonest improvements in the testimony of the prosecution witnesses. * The court also held that the acquittal of the co-accused by the Appellate Court was not suffering from any illegality or irregularity. "
